$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 10:46:14"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 10:46:08"
$wsZhCn.Range("K2").Value = "2016-08-22 10:46:33"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-22 10:46:14"
$wsDeDe.Range("K2").Value = "2016-08-22 10:46:40"
